# Commit: "Update 2p3. Added templates for formula student suspension,
# torque vectoring, four-wheel steering"
#
# Duplicate the existing Sedan_HambaLG sheet (keeping its full layout,
# column widths, styles and values) to create a new "FSAE_Achilles"
# template sheet placed right after it, then tweak the two cells that
# differ between the two templates: the "Instance" name shown in H3 and
# the rWheelCutout z-value in H6.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Copy the sheet as a new tab placed immediately after the source sheet.
$ws1.Copy($null, $ws1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "FSAE_Achilles"

# This new template is for the FSAE Achilles car, not the Sedan.
$ws2.Range("H3").Value = "FSAE_Achilles"

# Different rWheelCutout z value for this vehicle.
$ws2.Range("H6").Value = 0.25

# Make the new sheet the active tab, like it is right after being added.
$ws2.Activate()
